$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(12485,10962,10962,10962,10506,10413,10413,9502,9502,9502,9502,9036,9036,8989,8901,8901,8901,8901,8901,8789,8789,8487,8487,8487,8487,8487,8487,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8095,8056,8056,8056,8056,8056,8056,8056,8056,8056,8056,8056,8056,8056,8056,8056,8056,8056,8056,8056,7987,7987,7987,7987,7987,7987,7987,7987,7987,7987,7987,7987,7987,7987,7987,7987,7987,7987,7987,7987,7987,7987,7987,7987,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569,7569)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}
